$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C, shifting existing C:G to D:H
$ws.Columns("C:C").Insert()

# Match column B's width for the newly inserted blank column C (customWidth, no bestFit)
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Select G10 to match the resulting workbook state
$ws.Range("G10").Select()
